$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.690.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.59%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.099.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.48%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'229.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.34%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'61.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.88%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.81%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +2.44%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.92%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.412.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.57%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.52%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'22.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.99%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +2.28%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.70%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.106.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.85%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'38.622.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.51%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'71.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +3.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0836"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.56%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'226.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.40%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.04%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +4.42%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'170.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.93%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +7.70%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.11%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.80%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +6.28%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D35").Value = "'0.0607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.21%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.41%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.91%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.13%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.543.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.47%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'100.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.55%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.93%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.50%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +1.13%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.73%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'7.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +7.92%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +3.77%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.297.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.54%  "
$ws.Range("E51").Style = "Normal"
